# Generate Report for Archive
# - Update localization status text "Ready for handoff" -> "In Translation"
#   on every sheet / cell where it appears.
# - Shrink the now-narrower "Status" columns to match the new content width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-locale status columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
if ($wsOverview.Range("E2").Value2 -eq $oldStatus) { $wsOverview.Range("E2").Value = $newStatus }
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) { $wsOverview.Range("F2").Value = $newStatus }
if ($wsOverview.Range("E3").Value2 -eq $oldStatus) { $wsOverview.Range("E3").Value = $newStatus }
if ($wsOverview.Range("F3").Value2 -eq $oldStatus) { $wsOverview.Range("F3").Value = $newStatus }

# Columns E and F narrow now that "In Translation" is shorter than
# "Ready for handoff".
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Per-locale detail sheets: "Status" column C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    if ($ws.Range("C2").Value2 -eq $oldStatus) { $ws.Range("C2").Value = $newStatus }
    if ($ws.Range("C3").Value2 -eq $oldStatus) { $ws.Range("C3").Value = $newStatus }
    $ws.Columns.Item(3).ColumnWidth = 12.5
}
